# Kingdom 관련 모델 개편 (KingdomDeco, Structure, TIleMap) / ItemModel 추가
#
# In the "Packet" sheet of PlayerDetail.xlsx, the currency fields are
# renamed/repurposed in place:
#   Gold        -> Exp
#   AccGold     -> AccExp
#   StarCandy   -> Gold
#   AccStarCandy-> AccGold

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packet")

$ws.Range("A3").Value = "Exp"
$ws.Range("A4").Value = "AccExp"
$ws.Range("A5").Value = "Gold"
$ws.Range("A6").Value = "AccGold"
